$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B, C, D, E, F, G) for rows 2..11
$data = @{
    2  = @(-0.06578863049612835, 0.5873467334423524, 0.6636964401451462, 0.8146756656149404, 0.8199372453822267, 52)
    3  = @(-0.02523324468164586, 0.505465297153278,  0.4821656161413891, 0.6943814629880244, 0.7008277088516904, 51)
    4  = @(0.005608450684684903, 0.4847052646643364, 0.4358015571763177, 0.6601526771712115, 0.6668308405092866, 50)
    5  = @(-0.01284610207697986, 0.4921832099637495, 0.500471272416235,  0.7074399426214462, 0.7161102612933321, 41)
    6  = @(-0.09768167243710632, 0.5114466272893242, 0.4882265457959837, 0.6987320987302528, 0.7033071860062781, 31)
    7  = @(-0.1103781127069636,  0.5392641884755303, 0.5358859154982463, 0.7320422907853387, 0.7364823498096339, 29)
    8  = @(-0.09365692653440394, 0.5420077535482939, 0.5495523960454577, 0.7413180127620383, 0.7493864402156157, 27)
    9  = @(-0.0615712106401871,  0.5403710845059494, 0.5803707086519223, 0.7618206538627856, 0.7801358214643087, 19)
    10 = @(-0.01679119764512867, 0.4453590725437407, 0.4555235998480789, 0.6749248845968556, 0.704717858846984,  12)
    11 = @(0.3510802739590448,   0.6628255016558736, 0.9409865269452651, 0.970044600492815,  1.011020009805757,  5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
